# Auto-generated market-data refresh for Diabolos_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per scheduled runner snapshot.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 655.2105
$ws.Range("I33").Value = 326.84616
$ws.Range("J33").Value = 1366.6666
$ws.Range("K33").Value = 326.84616
$ws.Range("L33").Value = 1366.6666
$ws.Range("M33").Value = -97.84616
$ws.Range("N33").Value = -1824.6666
$ws.Range("H62").Value = 102072.625
$ws.Range("J62").Value = 127123
$ws.Range("L62").Value = 127123
$ws.Range("N62").Value = -128371
$ws.Range("H65").Value = 102072.625
$ws.Range("J65").Value = 127123
$ws.Range("L65").Value = 635615
$ws.Range("N65").Value = -641855
$ws.Range("H132").Value = 2192.7344
$ws.Range("I132").Value = 2199.0667
$ws.Range("K132").Value = 6597.2001
$ws.Range("M132").Value = -4067.2001
$ws.Range("H141").Value = 1668.6
$ws.Range("I141").Value = 1335.75
$ws.Range("K141").Value = 4007.25
$ws.Range("M141").Value = 1172.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1400.4286
$ws.Range("I2").Value = 1746
$ws.Range("J2").Value = 536.5
$ws.Range("K2").Value = 1746
$ws.Range("L2").Value = 536.5
$ws.Range("M2").Value = -1633
$ws.Range("N2").Value = -762.5
$ws.Range("H32").Value = 3073.4307
$ws.Range("J32").Value = 9752.77
$ws.Range("L32").Value = 9752.77
$ws.Range("N32").Value = -10326.77
$ws.Range("H97").Value = 443.61905
$ws.Range("I97").Value = 446.10526
$ws.Range("J97").Value = 420
$ws.Range("K97").Value = 446.10526
$ws.Range("L97").Value = 420
$ws.Range("M97").Value = 49.89474000000001
$ws.Range("N97").Value = -1412
$ws.Range("H116").Value = 1400.4286
$ws.Range("I116").Value = 1746
$ws.Range("J116").Value = 536.5
$ws.Range("K116").Value = 1746
$ws.Range("L116").Value = 536.5
$ws.Range("M116").Value = 548
$ws.Range("N116").Value = -5124.5
$ws.Range("H132").Value = 37039370
$ws.Range("I132").Value = 38463844
$ws.Range("K132").Value = 115391532
$ws.Range("M132").Value = -115389002

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1400.4286
$ws.Range("I3").Value = 1746
$ws.Range("J3").Value = 536.5
$ws.Range("K3").Value = 1746
$ws.Range("L3").Value = 536.5
$ws.Range("M3").Value = -1632
$ws.Range("N3").Value = -764.5
$ws.Range("H86").Value = 35718256
$ws.Range("J86").Value = 1241
$ws.Range("L86").Value = 1241
$ws.Range("N86").Value = -3487
$ws.Range("H89").Value = 35718256
$ws.Range("J89").Value = 1241
$ws.Range("L89").Value = 6205
$ws.Range("N89").Value = -17437
$ws.Range("H107").Value = 31250682
$ws.Range("I107").Value = 707.3570999999999
$ws.Range("K107").Value = 707.3570999999999
$ws.Range("M107").Value = 1212.6429
$ws.Range("H134").Value = 1689
$ws.Range("I134").Value = 1626.8
$ws.Range("K134").Value = 4880.4
$ws.Range("M134").Value = -2345.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3217927
$ws.Range("I6").Value = 3217927
$ws.Range("K6").Value = 3217927
$ws.Range("M6").Value = -3217814
$ws.Range("H31").Value = 3541.8035
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3541.8035
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3541.8035
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4131.8035
$ws.Range("H34").Value = 3541.8035
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3541.8035
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3541.8035
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3945.8035
$ws.Range("H107").Value = 2619.5715
$ws.Range("I107").Value = 1787.3636
$ws.Range("K107").Value = 1787.3636
$ws.Range("M107").Value = 132.6364000000001
$ws.Range("H132").Value = 2486.8147
$ws.Range("J132").Value = 5999
$ws.Range("L132").Value = 17997
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 1191.48
$ws.Range("I134").Value = 1035.1818
$ws.Range("J134").Value = 2337.6667
$ws.Range("K134").Value = 3105.5454
$ws.Range("L134").Value = 7013.000100000001
$ws.Range("M134").Value = -570.5454
$ws.Range("N134").Value = -12083.0001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10251899
$ws.Range("I4").Value = 11467787
$ws.Range("K4").Value = 34403361
$ws.Range("M4").Value = -34403249
$ws.Range("H7").Value = 10.6
$ws.Range("I7").Value = 10.666667
$ws.Range("K7").Value = 32.000001
$ws.Range("M7").Value = 79.999999
$ws.Range("H81").Value = 7259
$ws.Range("I81").Value = 1005.3333
$ws.Range("J81").Value = 8152.381
$ws.Range("K81").Value = 3015.9999
$ws.Range("L81").Value = 24457.143
$ws.Range("M81").Value = -1892.9999
$ws.Range("N81").Value = -26703.143
$ws.Range("H84").Value = 7259
$ws.Range("I84").Value = 1005.3333
$ws.Range("J84").Value = 8152.381
$ws.Range("K84").Value = 9047.9997
$ws.Range("L84").Value = 73371.429
$ws.Range("M84").Value = -3431.9997
$ws.Range("N84").Value = -84603.429
$ws.Range("H104").Value = 3211
$ws.Range("I104").Value = 1423
$ws.Range("J104").Value = 4999
$ws.Range("K104").Value = 4269
$ws.Range("L104").Value = 14997
$ws.Range("M104").Value = -1648
$ws.Range("N104").Value = -20239
$ws.Range("H121").Value = 120897.3
$ws.Range("J121").Value = 138390.38
$ws.Range("L121").Value = 415171.14
$ws.Range("N121").Value = -417791.14
$ws.Range("H129").Value = 1620.45
$ws.Range("I129").Value = 786.3333
$ws.Range("J129").Value = 2871.625
$ws.Range("K129").Value = 2358.9999
$ws.Range("L129").Value = 8614.875
$ws.Range("M129").Value = 2641.0001
$ws.Range("N129").Value = -18614.875
$ws.Range("H131").Value = 32224.572
$ws.Range("J131").Value = 44691.4
$ws.Range("L131").Value = 134074.2
$ws.Range("N131").Value = -144154.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 2335.4
$ws.Range("I55").Value = 2333.3333
$ws.Range("J55").Value = 2338.5
$ws.Range("K55").Value = 2333.3333
$ws.Range("L55").Value = 2338.5
$ws.Range("M55").Value = -2006.3333
$ws.Range("N55").Value = -2992.5
$ws.Range("H97").Value = 823.2778
$ws.Range("I97").Value = 600.4545000000001
$ws.Range("J97").Value = 1173.4286
$ws.Range("K97").Value = 600.4545000000001
$ws.Range("L97").Value = 1173.4286
$ws.Range("M97").Value = -104.4545000000001
$ws.Range("N97").Value = -2165.4286

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16669039
$ws.Range("I7").Value = 20835184
$ws.Range("K7").Value = 20835184
$ws.Range("M7").Value = -20835072
$ws.Range("H22").Value = 650
$ws.Range("J22").Value = 650
$ws.Range("L22").Value = 650
$ws.Range("N22").Value = -1240
$ws.Range("H27").Value = 650
$ws.Range("J27").Value = 650
$ws.Range("L27").Value = 650
$ws.Range("N27").Value = -864
$ws.Range("H40").Value = 4776.3335
$ws.Range("I40").Value = 4597.4
$ws.Range("K40").Value = 4597.4
$ws.Range("M40").Value = -4461.4
$ws.Range("H126").Value = 16669039
$ws.Range("I126").Value = 20835184
$ws.Range("K126").Value = 62505552
$ws.Range("M126").Value = -62503082
$ws.Range("H132").Value = 5217.9355
$ws.Range("I132").Value = 2715.389
$ws.Range("J132").Value = 8683
$ws.Range("K132").Value = 8146.167
$ws.Range("L132").Value = 26049
$ws.Range("M132").Value = -5616.167
$ws.Range("N132").Value = -31109

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8005913
$ws.Range("I81").Value = 3229.3076
$ws.Range("J81").Value = 16675487
$ws.Range("K81").Value = 6458.6152
$ws.Range("L81").Value = 33350974
$ws.Range("M81").Value = -5397.6152
$ws.Range("N81").Value = -33353096
$ws.Range("H84").Value = 8005913
$ws.Range("I84").Value = 3229.3076
$ws.Range("J84").Value = 16675487
$ws.Range("K84").Value = 32293.076
$ws.Range("L84").Value = 166754870
$ws.Range("M84").Value = -26989.076
$ws.Range("N84").Value = -166765478
$ws.Range("H122").Value = 2008.6
$ws.Range("I122").Value = 1886
$ws.Range("K122").Value = 5658
$ws.Range("M122").Value = -3208
$ws.Range("H132").Value = 3424.2354
$ws.Range("I132").Value = 3441.889
$ws.Range("K132").Value = 10325.667
$ws.Range("M132").Value = -7795.667000000001
$ws.Range("H136").Value = 4013.4443
$ws.Range("I136").Value = 1626.25
$ws.Range("J136").Value = 5923.2
$ws.Range("K136").Value = 4878.75
$ws.Range("L136").Value = 17769.6
$ws.Range("M136").Value = -2328.75
$ws.Range("N136").Value = -22869.6
